{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\"Mi programa\", \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$lastPara = $d.Paragraphs.Last\n$newRange = $lastPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Mi programa\"\n"}
